$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from header cell H1 into the new header cells I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set header values for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns (I = I0, J = IF)
$data = @(
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(7, 8),
    @(7, 8),
    @(1, 5),
    @(6, 6),
    @(6, 9),
    @(6, 7),
    @(6, 7),
    @(5, 6),
    @(5, 5),
    @(5, 7),
    @(6, 8),
    @(5, 6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
